# Auto-generated Excel COM-interop script applying scheduled-runner market data updates
# to the Valefor_Profits workbook (per-sheet leve-profit recalculation refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2759
$ws.Range("I62").Value = 2759
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2759
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2135
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 2759
$ws.Range("I65").Value = 2759
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13795
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -10675
$ws.Range("N65").ClearContents()

$ws.Range("H86").Value = 3643.3333
$ws.Range("I86").Value = 3637.5
$ws.Range("J86").Value = 3666.6667
$ws.Range("K86").Value = 3637.5
$ws.Range("L86").Value = 3666.6667
$ws.Range("M86").Value = -2514.5
$ws.Range("N86").Value = -5912.6667

$ws.Range("H89").Value = 3643.3333
$ws.Range("I89").Value = 3637.5
$ws.Range("J89").Value = 3666.6667
$ws.Range("K89").Value = 18187.5
$ws.Range("L89").Value = 18333.3335
$ws.Range("M89").Value = -12571.5
$ws.Range("N89").Value = -29565.3335

$ws.Range("H98").Value = 45469.76
$ws.Range("I98").Value = 54480.973
$ws.Range("J98").Value = 2666.5
$ws.Range("K98").Value = 54480.973
$ws.Range("L98").Value = 2666.5
$ws.Range("M98").Value = -52982.973
$ws.Range("N98").Value = -5662.5

$ws.Range("H106").Value = 11150
$ws.Range("I106").Value = 11150
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 11150
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -10519
$ws.Range("N106").ClearContents()

$ws.Range("H112").Value = 1124.8363
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 1133.3208
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 3399.9624
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -5615.9624

$ws.Range("H122").Value = 45469.76
$ws.Range("I122").Value = 54480.973
$ws.Range("J122").Value = 2666.5
$ws.Range("K122").Value = 163442.919
$ws.Range("L122").Value = 7999.5
$ws.Range("M122").Value = -160992.919
$ws.Range("N122").Value = -12899.5

$ws.Range("H129").Value = 1737.0667
$ws.Range("I129").Value = 2007
$ws.Range("J129").Value = 1428.5714
$ws.Range("K129").Value = 6021
$ws.Range("L129").Value = 4285.7142
$ws.Range("M129").Value = -1021
$ws.Range("N129").Value = -14285.7142

$ws.Range("H138").Value = 5067.0835
$ws.Range("I138").Value = 1897.4
$ws.Range("J138").Value = 6123.6445
$ws.Range("K138").Value = 5692.200000000001
$ws.Range("L138").Value = 18370.9335
$ws.Range("M138").Value = -552.2000000000007
$ws.Range("N138").Value = -28650.9335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25689.877
$ws.Range("I32").Value = 15114.462
$ws.Range("J32").Value = 67991.53999999999
$ws.Range("K32").Value = 15114.462
$ws.Range("L32").Value = 67991.53999999999
$ws.Range("M32").Value = -14827.462
$ws.Range("N32").Value = -68565.53999999999

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H74").Value = 10870715
$ws.Range("I74").Value = 12821514
$ws.Range("J74").Value = 1975.4286
$ws.Range("K74").Value = 12821514
$ws.Range("L74").Value = 1975.4286
$ws.Range("M74").Value = -12820640
$ws.Range("N74").Value = -3723.4286

$ws.Range("H77").Value = 10870715
$ws.Range("I77").Value = 12821514
$ws.Range("J77").Value = 1975.4286
$ws.Range("K77").Value = 64107570
$ws.Range("L77").Value = 9877.143
$ws.Range("M77").Value = -64103202
$ws.Range("N77").Value = -18613.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 932.5
$ws.Range("I22").Value = 597.5
$ws.Range("J22").Value = 1100
$ws.Range("K22").Value = 597.5
$ws.Range("L22").Value = 1100
$ws.Range("M22").Value = -247.5
$ws.Range("N22").Value = -1800

$ws.Range("H31").Value = 23813712
$ws.Range("I31").Value = 41668084
$ws.Range("J31").Value = 7886
$ws.Range("K31").Value = 41668084
$ws.Range("L31").Value = 7886
$ws.Range("M31").Value = -41667789
$ws.Range("N31").Value = -8476

$ws.Range("H34").Value = 23813712
$ws.Range("I34").Value = 41668084
$ws.Range("J34").Value = 7886
$ws.Range("K34").Value = 41668084
$ws.Range("L34").Value = 7886
$ws.Range("M34").Value = -41667882
$ws.Range("N34").Value = -8290

$ws.Range("H105").Value = 36355.766
$ws.Range("I105").Value = 47603.773
$ws.Range("J105").Value = 5423.75
$ws.Range("K105").Value = 47603.773
$ws.Range("L105").Value = 5423.75
$ws.Range("M105").Value = -45856.773
$ws.Range("N105").Value = -8917.75

$ws.Range("H107").Value = 1370.65
$ws.Range("I107").Value = 290.125
$ws.Range("J107").Value = 2091
$ws.Range("K107").Value = 290.125
$ws.Range("L107").Value = 2091
$ws.Range("M107").Value = 1629.875
$ws.Range("N107").Value = -5931

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 2791.2222
$ws.Range("I132").Value = 2103.3845
$ws.Range("J132").Value = 4579.6
$ws.Range("K132").Value = 6310.1535
$ws.Range("L132").Value = 13738.8
$ws.Range("M132").Value = -3780.1535
$ws.Range("N132").Value = -18798.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 39
$ws.Range("I12").Value = 19.5
$ws.Range("J12").Value = 44.086956
$ws.Range("K12").Value = 58.5
$ws.Range("L12").Value = 132.260868
$ws.Range("M12").Value = 114.5
$ws.Range("N12").Value = -478.260868

$ws.Range("H114").Value = 1636.3636
$ws.Range("J114").Value = 5150
$ws.Range("L114").Value = 15450
$ws.Range("N114").Value = -21958

$ws.Range("H121").Value = 40005612
$ws.Range("J121").Value = 41672410
$ws.Range("L121").Value = 125017230
$ws.Range("N121").Value = -125019850

$ws.Range("H122").Value = 850.65515
$ws.Range("I122").Value = 645.63635
$ws.Range("J122").Value = 1495
$ws.Range("K122").Value = 5810.72715
$ws.Range("L122").Value = 13455
$ws.Range("M122").Value = -3360.72715
$ws.Range("N122").Value = -18355

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 48754.273
$ws.Range("I80").Value = 2270
$ws.Range("J80").Value = 87491.164
$ws.Range("K80").Value = 2270
$ws.Range("L80").Value = 87491.164
$ws.Range("M80").Value = -1272
$ws.Range("N80").Value = -89487.164

$ws.Range("H83").Value = 48754.273
$ws.Range("I83").Value = 2270
$ws.Range("J83").Value = 87491.164
$ws.Range("K83").Value = 11350
$ws.Range("L83").Value = 437455.82
$ws.Range("M83").Value = -6358
$ws.Range("N83").Value = -447439.82

$ws.Range("H122").Value = 1340
$ws.Range("I122").Value = 1340
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4020
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1570
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1242.4
$ws.Range("I22").Value = 2001
$ws.Range("J22").Value = 1052.75
$ws.Range("K22").Value = 2001
$ws.Range("L22").Value = 1052.75
$ws.Range("M22").Value = -1706
$ws.Range("N22").Value = -1642.75

$ws.Range("H27").Value = 1242.4
$ws.Range("I27").Value = 2001
$ws.Range("J27").Value = 1052.75
$ws.Range("K27").Value = 2001
$ws.Range("L27").Value = 1052.75
$ws.Range("M27").Value = -1894
$ws.Range("N27").Value = -1266.75

$ws.Range("H40").Value = 3125.5
$ws.Range("I40").Value = 2667.3333
$ws.Range("J40").Value = 4500
$ws.Range("K40").Value = 2667.3333
$ws.Range("L40").Value = 4500
$ws.Range("M40").Value = -2531.3333
$ws.Range("N40").Value = -4772

$ws.Range("H82").Value = 2437
$ws.Range("I82").Value = 1750.6666
$ws.Range("K82").Value = 1750.6666
$ws.Range("M82").Value = -1389.6666

$ws.Range("H85").Value = 2437
$ws.Range("I85").Value = 1750.6666
$ws.Range("K85").Value = 1750.6666
$ws.Range("M85").Value = -502.6666

$ws.Range("H122").Value = 6707.92
$ws.Range("I122").Value = 6779.0835
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 20337.2505
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -17887.2505
$ws.Range("N122").Value = -19900

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 3966.3333
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3966.3333
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 11898.9999
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -16958.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 10000
$ws.Range("I17").Value = 10000
$ws.Range("K17").Value = 10000
$ws.Range("M17").Value = -9828

$ws.Range("H109").Value = 34171.43
$ws.Range("J109").Value = 34166.668
$ws.Range("L109").Value = 34166.668
$ws.Range("N109").Value = -36940.668

$ws.Range("H136").Value = 9503.588
$ws.Range("I136").Value = 3871.4285
$ws.Range("J136").Value = 13446.1
$ws.Range("K136").Value = 11614.2855
$ws.Range("L136").Value = 40338.3
$ws.Range("M136").Value = -9064.2855
$ws.Range("N136").Value = -45438.3
